# Append JPY/USD FX daily OHLCV rows for 2023-09-01 .. 2023-10-24
# (rows 813-850) to the end of the existing data table, matching the
# formatting of the adjacent existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 813
$lastNewRow  = 850

# Copy the format (border / bold font / date number format / alignment)
# of the last existing data row down across the full new range so the
# new rows look identical to the ones above them.
$ws.Range("A812:F812").Copy()
$ws.Range("A${firstNewRow}:F${lastNewRow}").PasteSpecial(-4122)

# date-serial, open, high, low, close, volume
$data = @(
    @(813, 45170, 0.00687, 0.006921, 0.006836, 0.006836, 0),
    @(814, 45173, 0.006836, 0.006847, 0.006825, 0.006825, 0),
    @(815, 45174, 0.006826, 0.006826, 0.006767, 0.006769, 0),
    @(816, 45175, 0.006769, 0.006801, 0.006767, 0.006771, 0),
    @(817, 45176, 0.006772, 0.0068, 0.006764, 0.006786, 0),
    @(818, 45177, 0.006787, 0.006818, 0.006764, 0.006764, 0),
    @(819, 45180, 0.006764, 0.006852, 0.006764, 0.00682, 0),
    @(820, 45181, 0.006821, 0.006826, 0.006793, 0.006796, 0),
    @(821, 45182, 0.006795, 0.006801, 0.006774, 0.006783, 0),
    @(822, 45183, 0.006783, 0.006801, 0.006778, 0.006779, 0),
    @(823, 45184, 0.00678, 0.006786, 0.00676, 0.006761, 0),
    @(824, 45187, 0.006761, 0.006778, 0.006761, 0.006773, 0),
    @(825, 45188, 0.006774, 0.006778, 0.006762, 0.006764, 0),
    @(826, 45189, 0.006763, 0.00678, 0.006741, 0.006745, 0),
    @(827, 45190, 0.006744, 0.006787, 0.006737, 0.006774, 0),
    @(828, 45191, 0.006774, 0.006777, 0.006737, 0.006737, 0),
    @(829, 45194, 0.006737, 0.006743, 0.006714, 0.006716, 0),
    @(830, 45195, 0.006716, 0.006723, 0.006704, 0.006708, 0),
    @(831, 45196, 0.006709, 0.006716, 0.00668, 0.006686, 0),
    @(832, 45197, 0.006687, 0.006703, 0.006687, 0.006696, 0),
    @(833, 45198, 0.006696, 0.006732, 0.00669, 0.006692, 0),
    @(834, 45201, 0.006692, 0.006694, 0.006667, 0.006672, 0),
    @(835, 45202, 0.006672, 0.006757, 0.006661, 0.006707, 0),
    @(836, 45203, 0.006704, 0.006722, 0.006698, 0.006708, 0),
    @(837, 45204, 0.006708, 0.006742, 0.006708, 0.006732, 0),
    @(838, 45205, 0.006733, 0.006739, 0.00669, 0.006697, 0),
    @(839, 45208, 0.006697, 0.006736, 0.006679, 0.006731, 0),
    @(840, 45209, 0.006731, 0.006748, 0.006708, 0.006722, 0),
    @(841, 45210, 0.006723, 0.006736, 0.006698, 0.006705, 0),
    @(842, 45211, 0.006705, 0.006713, 0.006674, 0.006674, 0),
    @(843, 45212, 0.006675, 0.00669, 0.006675, 0.006683, 0),
    @(844, 45215, 0.006683, 0.006695, 0.006678, 0.006686, 0),
    @(845, 45216, 0.006686, 0.006695, 0.006674, 0.006674, 0),
    @(846, 45217, 0.006674, 0.006686, 0.006669, 0.006673, 0),
    @(847, 45218, 0.006673, 0.00668, 0.00667, 0.006675, 0),
    @(848, 45219, 0.006673, 0.006676, 0.006668, 0.006671, 0),
    @(849, 45222, 0.006671, 0.006686, 0.006666, 0.006677, 0),
    @(850, 45223, 0.0066776, 0.006688, 0.0066769, 0.006686, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
